$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 74.02
$ws.Range("C2").Value = "The applicant possesses strong skills in TensorFlow, JavaScript, Scikit-learn, and Python which are crucial for the job role. Additionally, the applicant has experience in Django, ReactJS, and image processing, showcasing versatility and adaptability. However, the applicant lacks experience in PyTorch, and this could be a skill gap for the role. Nevertheless, the applicant's proficiency in Computer Vision and data analysis complements the job requirements, making the applicant a suitable candidate overall."
$ws.Range("E2").Value = "The applicant has exhibited a proactive and adaptable personality during the interview process, showcasing a willingness to embrace new opportunities. Their responses also indicate an open-minded approach and a high level of flexibility, which align with the soft skills required by the company. Overall, the applicant demonstrates a strong willingness to work in Japan, making them a good fit for the company's international expansion plans."

$ws.Range("B3").Value = 89.42
$ws.Range("C3").Value = "The applicant has experience with ReactJS, Flutter, and NodeJS which are essential for this job role. However, the applicant lacks experience with WebRTC, Socket.IO, and Docker, which are required by the company. Despite having some relevant skills, the lack of experience in these key areas makes the applicant less suited for this job position."
$ws.Range("E3").Value = "The applicant has demonstrated strong communication skills and adaptability, making them suitable for the company's soft skill requirements. Additionally, the applicant has shown an eagerness to experience new cultures, which aligns with the company's focus on diversity and global perspectives."

$ws.Range("B4").Value = 77.25
$ws.Range("C4").Value = "The applicant possesses skills in PyTorch, TensorFlow, Django Rest Framework, Scikit-learn, and Keras, all of which are directly relevant to the job description provided by the company. The applicant also has experience in computer vision and image processing, which are additional valuable skills for the role. However, the applicant lacks experience in ReactJS and Python, which are among the required skills by the company. As a result, though the applicant has a strong background in the required technologies, the lack of experience in ReactJS and Python may impact the overall suitability for the role."
$ws.Range("E4").Value = "The applicant has demonstrated good adaptability and openness to new experiences, which are essential qualities for working in a new environment such as Japan. Additionally, the applicant's responses indicated a high level of cultural awareness and respect, aligning well with the interpersonal skills sought by the company. Therefore, based on the applicant's personality traits and willingness to work in a new cultural setting, a score of 4 is justified."

$ws.Range("B5").Value = 82.91
$ws.Range("C5").Value = "The applicant has relevant skills in API integration, JavaScript, HTML, and CSS, which align with the required skills set by the company. However, the applicant lacks experience in specific technologies mentioned in the job description, such as Python, Java, and AI/ML. While the projects demonstrate proficiency in web development, the absence of experience in the required technologies may limit the applicant's suitability for the role."
$ws.Range("E5").Value = "The applicant exhibits a willingness to adapt and embrace new challenges, with a positive attitude towards relocation. Their responses reflect a strong sense of adaptability and open-mindedness, which are essential traits for thriving in a new environment. Moreover, the applicant demonstrates a keen interest in embracing the cultural aspects of Japan, indicating a genuine enthusiasm for the opportunity."
